$wb = $excel.ActiveWorkbook

# --- Sheet: paper ---
$ws = $wb.Worksheets.Item("paper")
$rows = @(
    @{Row=2; B=3.70179780381; C=10.296678571059; D=0},
    @{Row=3; B=3.912801029183; C=21.880200560139; D=0},
    @{Row=4; B=7.885903670996; C=24.868462576785; D=0},
    @{Row=5; B=2.727868101063; C=3.763599152711; D=0},
    @{Row=6; B=3.431908146809; C=6.635476847929; D=0},
    @{Row=7; B=14.994269705273; C=4.559949701541; D=0},
    @{Row=8; B=3.033003614842; C=3.650595717278; D=0},
    @{Row=9; B=2.008297157096; C=5.537623371892; D=0},
    @{Row=10; B=3.828234237744; C=6.632690200563; D=0},
    @{Row=11; B=2.203260542289; C=2.016655523388; D=0},
    @{Row=12; B=3.255684611193; C=5.836446991391; D=0},
    @{Row=13; B=0; C=0; D=0},
    @{Row=14; B=16.900208051391; C=16.27493192542; D=0},
    @{Row=15; B=3.935004225548; C=8.783221028129001; D=0},
    @{Row=16; B=3.292797835025; C=5.722414507206; D=0},
    @{Row=17; B=3.264294606093; C=10.20499056547; D=0},
    @{Row=18; B=2.877348880723; C=6.584127129136; D=0},
    @{Row=19; B=5.086326162223; C=17.831986095087; D=0},
    @{Row=20; B=5.602795956783; C=9.066451183926; D=0},
    @{Row=21; B=3.046888646976; C=7.446027371944; D=0},
    @{Row=22; B=4.921403167672; C=19.876389962204; D=0},
    @{Row=23; B=6.262872904277; C=17.76790957815; D=0},
    @{Row=24; B=6.209678077252; C=14.52685264113; D=0},
    @{Row=25; B=14.188280019771; C=11.898159210155; D=0},
    @{Row=26; B=11.521456568215; C=6.408607878133; D=0},
    @{Row=27; B=2.994690265487; C=7.470796460177; D=0},
    @{Row=28; B=2.994690265487; C=7.470796460177; D=0},
    @{Row=29; B=0.561635567418; C=2.389627444573; D=0},
    @{Row=30; B=2.994690265487; C=7.470796460177; D=0},
    @{Row=31; B=3.819869337975; C=6.926106787373; D=0},
    @{Row=32; B=2.994690265487; C=7.470796460177; D=0},
    @{Row=33; B=0; C=0; D=0},
    @{Row=34; B=2.715597004453; C=0; D=0},
    @{Row=35; B=20.963610949203; C=13.079944754009; D=0}
)
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# --- Sheet: glass ---
$ws = $wb.Worksheets.Item("glass")
$rows = @(
    @{Row=2; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=3; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=4; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=5; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=6; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=7; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=8; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=9; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=10; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=11; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=12; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=13; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=14; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=15; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=16; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=17; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=18; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=19; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=20; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=21; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=22; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=23; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=24; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=25; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=26; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=27; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=28; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=29; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=30; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=31; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=32; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=33; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=34; B=1.389896080003; C=6.133162642233; D=0},
    @{Row=35; B=1.389896080003; C=6.133162642233; D=0}
)
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
